# Fix the "Time in session" and "Public bills enacted into law" rows so the
# previously-merged Senate/House values live in their own cells, and correct
# the "Special reporcs" typo to "Special reports".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 - fix typo "Special reporcs " -> "Special reports "
$ws.Range("A28").Value = "Special reports "

# Row 10 - "Public bills enacted into law": split "28 60" (which lived only
# in C10 as text) into B10 (Senate) = 28 and C10 (House) = 60, as numbers
$ws.Range("B10").Value = 28
$ws.Range("C10").Value = 60

# Row 7 - "Time in session": split "1,839 hrs., 10' 1,525 hrs., 25'" (which
# lived only in C7) into B7 (Senate) = "1,839 hrs., 10' " and
# C7 (House) = "1,525 hrs., 25'"
$ws.Range("C7").Value = "1,525 hrs., 25'"
$ws.Range("B7").Value = "1,839 hrs., 10' "
